$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F4 and F5
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1372
$wsExhibit.Range("F5").Value = 657

# Sheet "全部类型" (sheet4): update F4 and F6
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1372
$wsAll.Range("F6").Value = 657
